$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '36.897.86'
Set-TextValue "E2" '  -1.66%  '
Set-TextValue "D3" '2.011.75'
Set-TextValue "E3" '  -2.68%  '
Set-TextValue "D4" '1.01'
Set-TextValue "E4" '  +0.51%  '
Set-TextValue "D5" '225.18'
Set-TextValue "E5" '  -2.85%  '
Set-TextValue "D6" '0.605'
Set-TextValue "E6" '  -3.66%  '
Set-TextValue "D8" '54.52'
Set-TextValue "E8" '  -4.84%  '
Set-TextValue "D9" '0.378'
Set-TextValue "E9" '  -2.53%  '
Set-TextValue "D10" '0.0781'
Set-TextValue "E10" '  +0.54%  '
Set-TextValue "E11" '  -4.13%  '
Set-TextValue "D12" '2.315.93'
Set-TextValue "E12" '  -2.45%  '
Set-TextValue "D13" '14.19'
Set-TextValue "E13" '  -4.41%  '
Set-TextValue "D14" '20.26'
Set-TextValue "E14" '  -2.62%  '
Set-TextValue "D15" '0.737'
Set-TextValue "E15" '  -3.06%  '
Set-TextValue "D16" '5.11'
Set-TextValue "E16" '  -3.60%  '
Set-TextValue "D17" '2.011.75'
Set-TextValue "E17" '  -2.65%  '
Set-TextValue "D18" '36.816.00'
Set-TextValue "E18" '  -1.74%  '
Set-TextValue "D19" '6.19'
Set-TextValue "E19" '  +4.83%  '
Set-TextValue "D20" '68.64'
Set-TextValue "E20" '  -2.46%  '
Set-TextValue "D21" '0.0₃0817'
Set-TextValue "E21" '  -1.25%  '
Set-TextValue "D22" '225.32'
Set-TextValue "E22" '  -1.05%  '
Set-TextValue "D23" '0.999'
Set-TextValue "E23" '  -0.12%  '
Set-TextValue "E24" '  +3.07%  '
Set-TextValue "D25" '2.17'
Set-TextValue "E25" '  -8.02%  '
Set-TextValue "D26" '165.15'
Set-TextValue "E26" '  -2.24%  '
Set-TextValue "D27" '9.16'
Set-TextValue "E27" '  -4.39%  '
Set-TextValue "E28" '  -5.32%  '
Set-TextValue "D29" '18.62'
Set-TextValue "E29" '  -3.76%  '
Set-TextValue "D30" '1.32'
Set-TextValue "E30" '  -3.10%  '
Set-TextValue "D31" '0.116'
Set-TextValue "E31" '  -4.59%  '
Set-TextValue "D32" '4.49'
Set-TextValue "E32" '  -1.81%  '
Set-TextValue "D33" '0.0614'
Set-TextValue "D34" '4.40'
Set-TextValue "E34" '  -4.04%  '
Set-TextValue "D35" '2.33'
Set-TextValue "E35" '  -5.17%  '
Set-TextValue "E36" '  +1.49%  '
Set-TextValue "E37" '  +0.26%  '
Set-TextValue "D38" '3.12'
Set-TextValue "E38" '  -5.54%  '
Set-TextValue "D39" '5.27'
Set-TextValue "E39" '  +0.15%  '
Set-TextValue "B40" 'Maker'
Set-TextValue "C40" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D40" '1.479.05'
Set-TextValue "E40" '  +0.56%  '
Set-TextValue "D41" '0.0216'
Set-TextValue "E41" '  -5.48%  '
Set-TextValue "B42" 'InjectiveProtocol'
Set-TextValue "C42" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D42" '16.92'
Set-TextValue "E42" '  +2.89%  '
Set-TextValue "D43" '94.85'
Set-TextValue "E43" '  -4.32%  '
Set-TextValue "D44" '0.0925'
Set-TextValue "E44" '  -3.19%  '
Set-TextValue "E45" '  -5.23%  '
Set-TextValue "E46" '  -6.25%  '
Set-TextValue "D47" '7.29'
Set-TextValue "E47" '  +1.00%  '
Set-TextValue "D48" '0.999'
Set-TextValue "E48" '  -3.65%  '
Set-TextValue "D49" '2.91'
Set-TextValue "E49" '  -0.98%  '
Set-TextValue "D50" '2.205.50'
Set-TextValue "E50" '  -2.34%  '
Set-TextValue "D51" '3.58'
Set-TextValue "E51" '  -10.28%  '
